$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.089.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "'1.873.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.03%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'313.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").Value = "'0.5047"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.65%  "
$ws.Range("D8").Value = "'0.3836"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.20%  "
$ws.Range("D9").Value = "'0.08668"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.80%  "
$ws.Range("E10").Value = "  -2.15%  "
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("D12").Value = "'6.330"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("D14").Value = "'1.872.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").Value = "'1.003"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").Value = "'7.167"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.00%  "
$ws.Range("D17").Value = "'0.00001103"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.64%  "
$ws.Range("D18").Value = "'90.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.59%  "
$ws.Range("D19").Value = "'0.06629"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("D20").Value = "'18.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").Value = "'6.104"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.96%  "
$ws.Range("D23").Value = "'28.128.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.49%  "
$ws.Range("D24").Value = "'11.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'2.261"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.51%  "
$ws.Range("D26").Value = "'2.567"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.73%  "
$ws.Range("D27").Value = "'2.087.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.88%  "
$ws.Range("E28").Value = "  -1.72%  "
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").Value = "'126.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("E31").Value = "  -2.41%  "
$ws.Range("E32").Value = "  -3.60%  "
$ws.Range("D33").Value = "'5.592"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'3.600"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("D35").Value = "'9.664"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.50%  "
$ws.Range("D36").Value = "'0.02460"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.56%  "
$ws.Range("D37").Value = "'0.06596"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.14%  "
$ws.Range("D38").Value = "'0.2173"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.17%  "
$ws.Range("D39").Value = "'1.205"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'1.244"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.02%  "
$ws.Range("E41").Value = "  +0.80%  "
$ws.Range("D42").Value = "'0.6367"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.19%  "
$ws.Range("D43").Value = "'4.897"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.83%  "
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").Value = "'13.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("D46").Value = "'0.5979"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.90%  "
$ws.Range("D47").Value = "'1.280"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("D48").Value = "'3.677"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("D49").Value = "'1.230"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.81%  "
$ws.Range("D50").Value = "'1.988"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("D51").Value = "'121.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.40%  "
